# Updated burn down chart
# Fill in "Actual" values (column C) for rows 13 through 19 with 6,
# matching the continued flat actual-progress line on the burndown chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C13").Value = 6
$ws.Range("C14").Value = 6
$ws.Range("C15").Value = 6
$ws.Range("C16").Value = 6
$ws.Range("C17").Value = 6
$ws.Range("C18").Value = 6
$ws.Range("C19").Value = 6

# Restore the view to scroll back up and adjust the zoom level, as seen
# after the user reviewed the updated chart.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.Zoom = 104

$ws.Range("A37").Select()
